# BOM MAIN BOARD.xlsx — replace two 4-pin JST SH connectors (J1 + J2,
# both described by a single "ESC, TELEM" BOM line) with a single 6-pin
# JST SH connector used for J2 only, and append the new crimp/cable/
# housing line items that go with it (rows 33-37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lrm = [char]0x200E
$mpn6pin = $lrm + "BM06B-SRSS-TB(LF)(SN)" + $lrm

# --- Row 11: J1/J2 4-pin connector line -> J2-only 6-pin connector line ---
# (Write order matches the new shared-string insertion order: Parts,
#  then Value/Device/MPN, then Package. Description text is unchanged.)
$ws.Range("A11").Value = 1
$ws.Range("E11").Value = "J2"
$ws.Range("B11").Value = $mpn6pin
$ws.Range("C11").Value = $mpn6pin
$ws.Range("G11").Value = $mpn6pin
$ws.Range("D11").Value = "JST SH 1MM 6 PIN"

# --- Rows 33-37: update/add the crimp/cable/housing line items ---
# MPN (col G) column is filled in first, top to bottom, then the
# Description (col F) column, top to bottom - matching the order the
# new shared strings were authored in.
$ws.Range("G33").Value = $lrm + "SSH-003T-P0.2"
$ws.Range("G34").Value = "AWG28-08/F-1/300"
$ws.Range("G35").Value = "SHR-06V-S"
$ws.Range("G36").Value = $lrm + "M20-1060400"
$ws.Range("G37").Value = "M20-1160042" + $lrm

$ws.Range("F33").Value = "CONN SOCKET 28-32AWG CRIMP TIN"
$ws.Range("F34").Value = "CBL RIBN 8COND 0.039 MULTI 5'"
$ws.Range("F35").Value = "CONN HOUSING SH 6POS 1MM WHITE"
$ws.Range("F36").Value = "CONN RCPT HSG 4POS 2.54MM"
$ws.Range("F37").Value = "CONN SOCKET 22-30AWG CRIMP GOLD"

$ws.Range("A33").Value = 6
$ws.Range("A34").Value = 1
$ws.Range("A35").Value = 1
$ws.Range("A36").Value = 2
$ws.Range("A37").Value = 8

# Rows 33-37 all carry a tall (wrapped-description) row height.
$ws.Rows.Item(34).RowHeight = 28.8
$ws.Rows.Item(35).RowHeight = 28.8
$ws.Rows.Item(36).RowHeight = 28.8
$ws.Rows.Item(37).RowHeight = 28.8

# Scroll position / selection left by the author after the edit.
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("E41").Select()
